$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the new bottom rows (16, 17) inherit the same formatting as the
# existing indexed rows (column A uses a bold/bordered/centered style).
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

# Shift existing "extr" rows (old rows 8-15) down by two rows (new rows 10-17),
# moving from the bottom up so we don't clobber data before reading it.
for ($old = 15; $old -ge 8; $old--) {
    $new = $old + 2
    $ws.Range("A$new").Value2 = $ws.Range("A$old").Value2
    $ws.Range("B$new").Value2 = $ws.Range("B$old").Value2
    $ws.Range("C$new").Value2 = $ws.Range("C$old").Value2
    $ws.Range("D$new").Value2 = $ws.Range("D$old").Value2
    $ws.Range("E$new").Value2 = $ws.Range("E$old").Value2
}

# New row 8: line7
$ws.Range("A8").Value2 = 6
$ws.Range("B8").Value2 = "line7"
$ws.Range("C8").Value2 = 14
$ws.Range("D8").Value2 = 11
$ws.Range("E8").Value2 = $true

# New row 9: line8
$ws.Range("A9").Value2 = 7
$ws.Range("B9").Value2 = "line8"
$ws.Range("C9").Value2 = 16
$ws.Range("D9").Value2 = 9
$ws.Range("E9").Value2 = $false

# Updated A-column index values for shifted extr rows (now rows 10-17)
$ws.Range("A10").Value2 = 8
$ws.Range("A11").Value2 = 9
$ws.Range("A12").Value2 = 10
$ws.Range("A13").Value2 = 11
$ws.Range("A14").Value2 = 12
$ws.Range("A15").Value2 = 13
$ws.Range("A16").Value2 = 14
$ws.Range("A17").Value2 = 15

# extr1 and extr2 (now rows 10 and 11) flip in_service False -> True
$ws.Range("E10").Value2 = $true
$ws.Range("E11").Value2 = $true
